# Update the two-digit division worksheet numbers (row,col -> new value)
# Uses direct Range.Text assignment (scoped per cell) instead of Find/Replace
# because duplicate operand strings exist in several cells and a document-wide
# Find could otherwise land on the wrong occurrence.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "58÷2="
$t.Cell(1,2).Range.Text = "95÷4="
$t.Cell(1,3).Range.Text = "72÷7="
$t.Cell(1,4).Range.Text = "62÷8="
$t.Cell(1,5).Range.Text = "36÷4="

$t.Cell(5,1).Range.Text = "36÷3="
$t.Cell(5,2).Range.Text = "76÷9="
$t.Cell(5,3).Range.Text = "41÷5="
$t.Cell(5,4).Range.Text = "17÷9="
$t.Cell(5,5).Range.Text = "81÷6="

$t.Cell(9,1).Range.Text = "34÷4="
$t.Cell(9,2).Range.Text = "82÷8="
$t.Cell(9,3).Range.Text = "13÷4="
$t.Cell(9,4).Range.Text = "89÷4="
$t.Cell(9,5).Range.Text = "49÷2="

$t.Cell(13,1).Range.Text = "93÷2="
$t.Cell(13,2).Range.Text = "35÷7="
$t.Cell(13,3).Range.Text = "25÷6="
$t.Cell(13,4).Range.Text = "17÷4="
$t.Cell(13,5).Range.Text = "93÷9="

$t.Cell(17,1).Range.Text = "96÷6="
$t.Cell(17,2).Range.Text = "31÷4="
$t.Cell(17,3).Range.Text = "47÷5="
$t.Cell(17,4).Range.Text = "42÷6="
$t.Cell(17,5).Range.Text = "19÷3="

